$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right marking 4 -> 5, Wrong marking -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 84 -> 105, Wrong total -2 -> -2.4, Max summary "82/112" -> "102.6/140"
$ws.Range("B12").Value = 105
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "102.6/140"
